$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "62.782.41"
Set-TextValue $ws.Range("E2") "  +0.15%  "
Set-TextValue $ws.Range("D3") "3.450.67"
Set-TextValue $ws.Range("E3") "  -0.26%  "
Set-TextValue $ws.Range("D4") "0.999"
Set-TextValue $ws.Range("E4") "  -0.04%  "
Set-TextValue $ws.Range("D5") "578.98"
Set-TextValue $ws.Range("E5") "  +0.21%  "
Set-TextValue $ws.Range("D6") "149.56"
Set-TextValue $ws.Range("E6") "  +2.37%  "
Set-TextValue $ws.Range("D7") "0.999"
Set-TextValue $ws.Range("E7") "  -0.16%  "
Set-TextValue $ws.Range("D8") "0.486"
Set-TextValue $ws.Range("E8") "  +0.69%  "
Set-TextValue $ws.Range("D9") "8.02"
Set-TextValue $ws.Range("E9") "  +5.35%  "
Set-TextValue $ws.Range("E10") "  -0.41%  "
Set-TextValue $ws.Range("D11") "0.416"
Set-TextValue $ws.Range("E11") "  +4.06%  "
Set-TextValue $ws.Range("D12") "4.042.62"
Set-TextValue $ws.Range("E12") "  -0.16%  "
Set-TextValue $ws.Range("E13") "  -0.13%  "
Set-TextValue $ws.Range("D14") "28.35"
Set-TextValue $ws.Range("E14") "  -4.67%  "
Set-TextValue $ws.Range("D15") "3.450.14"
Set-TextValue $ws.Range("E15") "  -0.57%  "
Set-TextValue $ws.Range("D16") "0.0000173"
Set-TextValue $ws.Range("E16") "  +1.47%  "
Set-TextValue $ws.Range("D17") "62.820.88"
Set-TextValue $ws.Range("E17") "  +0.07%  "
Set-TextValue $ws.Range("D18") "6.41"
Set-TextValue $ws.Range("E18") "  +0.97%  "
Set-TextValue $ws.Range("D19") "14.64"
Set-TextValue $ws.Range("E19") "  +1.80%  "
Set-TextValue $ws.Range("D20") "9.02"
Set-TextValue $ws.Range("E20") "  -2.06%  "
Set-TextValue $ws.Range("D21") "387.33"
Set-TextValue $ws.Range("E21") "  -0.24%  "
Set-TextValue $ws.Range("D22") "0.569"
Set-TextValue $ws.Range("E22") "  +1.09%  "
Set-TextValue $ws.Range("D23") "75.33"
Set-TextValue $ws.Range("E23") "  +0.55%  "
Set-TextValue $ws.Range("E24") "  +0.01%  "
Set-TextValue $ws.Range("D25") "3.585.09"
Set-TextValue $ws.Range("E25") "  -0.60%  "
Set-TextValue $ws.Range("E26") "  +1.64%  "
Set-TextValue $ws.Range("E27") "  +2.61%  "
Set-TextValue $ws.Range("D28") "7.71"
Set-TextValue $ws.Range("E28") "  +1.82%  "
Set-TextValue $ws.Range("D29") "1.00"
Set-TextValue $ws.Range("E29") "  +0.05%  "
Set-TextValue $ws.Range("D30") "8.02"
Set-TextValue $ws.Range("E30") "  -0.91%  "
Set-TextValue $ws.Range("E31") "  -0.87%  "
Set-TextValue $ws.Range("E32") "  -0.02%  "
Set-TextValue $ws.Range("E33") "  -3.52%  "
Set-TextValue $ws.Range("D34") "23.25"
Set-TextValue $ws.Range("E34") "  -1.99%  "
Set-TextValue $ws.Range("D35") "5.41"
Set-TextValue $ws.Range("E35") "  +2.63%  "
Set-TextValue $ws.Range("D36") "1.64"
Set-TextValue $ws.Range("E36") "  +4.81%  "
Set-TextValue $ws.Range("D37") "32.23"
Set-TextValue $ws.Range("E37") "  +2.81%  "
Set-TextValue $ws.Range("E38") "  -1.82%  "
Set-TextValue $ws.Range("D39") "169.25"
Set-TextValue $ws.Range("E39") "  -0.66%  "
Set-TextValue $ws.Range("D40") "3.484.68"
Set-TextValue $ws.Range("E40") "  -0.37%  "
Set-TextValue $ws.Range("D41") "0.0778"
Set-TextValue $ws.Range("E41") "  +1.36%  "
Set-TextValue $ws.Range("B42") "OKB"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D42") "42.85"
Set-TextValue $ws.Range("E42") "  +1.61%  "
Set-TextValue $ws.Range("B43") "Mantle"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D43") "0.786"
Set-TextValue $ws.Range("E43") "  -1.37%  "
Set-TextValue $ws.Range("E44") "  -2.14%  "
Set-TextValue $ws.Range("D45") "1.69"
Set-TextValue $ws.Range("E45") "  -0.72%  "
Set-TextValue $ws.Range("E46") "  -0.87%  "
Set-TextValue $ws.Range("D47") "2.569.16"
Set-TextValue $ws.Range("E47") "  -0.89%  "
Set-TextValue $ws.Range("D48") "6.92"
Set-TextValue $ws.Range("E48") "  +2.29%  "
Set-TextValue $ws.Range("D49") "2.23"
Set-TextValue $ws.Range("E49") "  +2.27%  "
Set-TextValue $ws.Range("D50") "22.59"
Set-TextValue $ws.Range("E50") "  -3.01%  "
Set-TextValue $ws.Range("E51") "  -0.02%  "
